$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value stored as serial 45189 (2023-09-20)
# for every data row from row 2 through row 66. Update it to 45190 (2023-09-21).
$ws.Range("C2:C66").Value = 45190
